# Add "2022-Q3" data: a new sheet inserted right after "总计" and before the
# existing "2022-Q1" sheet, plus a new summary row on "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet before the current "2022-Q1" sheet,
#    using the latter as the style/format template (identical header/row
#    layout across all the quarterly detail sheets).
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q1")
$newSheet = $wb.Worksheets.Add($template)
$newSheet.Name = "2022-Q3"

# Header row - identical text/format on every quarterly sheet, so copy it
# wholesale (values + formatting) from the template.
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4104)

# Row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'015061"
$newSheet.Range("C2").Value = "中信建投沪深300指数增强A"
$newSheet.Range("D2").Value = "'1.47"
$newSheet.Range("E2").Value = "'91.17"
$newSheet.Range("F2").Value = "'1.63"
$newSheet.Range("G2").Value = "'0.0240"
$newSheet.Range("H2").Value = 6

# Row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'015062"
$newSheet.Range("C3").Value = "中信建投沪深300指数增强C"
$newSheet.Range("D3").Value = "'1.14"
$newSheet.Range("E3").Value = "'91.17"
$newSheet.Range("F3").Value = "'1.63"
$newSheet.Range("G3").Value = "'0.0186"
$newSheet.Range("H3").Value = 6

# Copy the "index" column (A) formatting from the template sheet (bold,
# bordered, centred - same look as every other quarterly sheet).
$template.Range("A2:A3").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: a new 2022-Q3 summary row is inserted
#    at the top of the data (row 2) and every following row shifts down by
#    one, keeping a 0-based running index in column A.
# ---------------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 2
$totals.Range("D2").Value = 0.04

$totals.Range("B3").Value = "2022-Q1"
$totals.Range("C3").Value = 3
$totals.Range("D3").Value = 0.02

$totals.Range("B4").Value = "2021-Q4"
$totals.Range("C4").Value = 8
$totals.Range("D4").Value = 1.29

$totals.Range("B5").Value = "2021-Q3"
$totals.Range("C5").Value = 5
$totals.Range("D5").Value = 1.21

$totals.Range("B6").Value = "2021-Q2"
$totals.Range("C6").Value = 2
$totals.Range("D6").Value = 0.11

$totals.Range("B7").Value = "2021-Q1"
$totals.Range("C7").Value = 3
$totals.Range("D7").Value = 0.53

$totals.Range("B8").Value = "2020-Q4"
$totals.Range("C8").Value = 5
$totals.Range("D8").Value = 3.15

For ($i = 0; $i -le 6; $i++) {
    $totals.Cells.Item($i + 2, 1).Value = $i
}

# Match the "index" column style (bold / bordered / centred) already used in
# column A, and extend it down to the newly-added row 8.
$totals.Range("A2").Copy()
$totals.Range("A2:A8").PasteSpecial(-4122)
